$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> values that were in row 4 (date 44911 group, "Región de O'Higgins")
$ws.Range("D2").Value = 44911
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 220
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1000

# Row 3 -> values that were in row 5 (date 44911 group, "Región de O'Higgins")
$ws.Range("D3").Value = 44911
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 4000
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 800

# Row 4 -> values that were in row 2 (date 44915 group, "Provincia de Quillota")
$ws.Range("D4").Value = 44915
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 1200

# Row 5 -> values that were in row 3 (date 44915 group, "Provincia de Quillota")
$ws.Range("D5").Value = 44915
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("R5").Value = "Provincia de Quillota"
$ws.Range("S5").Value = 1000
